$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.547.33'
$ws.Range("E2").Value = '  -1.32%  '
$ws.Range("D3").Value = '''1.925.26'
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''243.72'
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").Value = '''0.9997'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = '''0.4714'
$ws.Range("E7").Value = '  -1.74%  '
$ws.Range("D8").Value = '''0.2885'
$ws.Range("E8").Value = '  -2.02%  '
$ws.Range("D9").Value = '''0.06792'
$ws.Range("E9").Value = '  +2.37%  '
$ws.Range("D10").Value = '''106.69'
$ws.Range("E10").Value = '  +6.06%  '
$ws.Range("D11").Value = '''18.41'
$ws.Range("E11").Value = '  -2.22%  '
$ws.Range("D12").Value = '''0.07766'
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("D13").Value = '''1.913.77'
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").Value = '''5.326'
$ws.Range("E14").Value = '  +3.20%  '
$ws.Range("D15").Value = '''0.6629'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '''293.81'
$ws.Range("E16").Value = '  -4.51%  '
$ws.Range("D17").Value = '''30.573.18'
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").Value = '''0.000007613'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = '''12.98'
$ws.Range("E19").Value = '  -1.86%  '
$ws.Range("D20").Value = '''1.000'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '''2.155.72'
$ws.Range("E21").Value = '  +1.68%  '
$ws.Range("D22").Value = '''5.362'
$ws.Range("E22").Value = '  +3.23%  '
$ws.Range("D23").Value = '''0.9989'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = '''6.241'
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("D25").Value = '''9.378'
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '''168.90'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '''21.36'
$ws.Range("E27").Value = '  +3.20%  '
$ws.Range("D28").Value = '''2.105'
$ws.Range("E28").Value = '  +6.30%  '
$ws.Range("D29").Value = '''0.1076'
$ws.Range("E29").Value = '  -5.02%  '
$ws.Range("D30").Value = '''1.367'
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("D31").Value = '''4.201'
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").Value = '''4.009'
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("D33").Value = '''0.05053'
$ws.Range("E33").Value = '  -1.32%  '
$ws.Range("D34").Value = '''0.7426'
$ws.Range("E34").Value = '  -0.88%  '
$ws.Range("E35").Value = '  -1.22%  '
$ws.Range("D36").Value = '''0.02108'
$ws.Range("E36").Value = '  +5.62%  '
$ws.Range("D37").Value = '''2.727'
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").Value = '''2.688'
$ws.Range("E38").Value = '  -0.77%  '
$ws.Range("D39").Value = '''2.082'
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("D40").Value = '''110.60'
$ws.Range("E40").Value = '  +1.50%  '
$ws.Range("D41").Value = '''0.8790'
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("D42").Value = '''5.937'
$ws.Range("E42").Value = '  +4.51%  '
$ws.Range("D43").Value = '''0.4291'
$ws.Range("E43").Value = '  +1.09%  '
$ws.Range("D44").Value = '''1.000'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").Value = '''67.94'
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("D46").Value = '''50.05'
$ws.Range("E46").Value = '  +17.18%  '
$ws.Range("D47").Value = '''7.242'
$ws.Range("E47").Value = '  -1.97%  '
$ws.Range("D48").Value = '''9.338'
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("D49").Value = '''0.1223'
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").Value = '''35.21'
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("D51").Value = '''0.2470'
$ws.Range("E51").Value = '  +8.66%  '
